$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '60.684.17'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.23%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.641.21'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +1.24%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '575.34'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.64%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '143.87'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.10%  '
$ws.Range('E7').Value = '  +0.26%  '
$ws.Range('E8').Value = '  -0.45%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '6.57'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.47%  '
$ws.Range('E10').Value = '  +0.09%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.379'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +2.48%  '
$ws.Range('E12').Value = '  -0.62%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.116.28'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.46%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '26.14'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +11.66%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '60.752.27'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.07%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000143'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.30%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.660.34'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +1.51%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '11.53'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +2.30%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.72'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.99%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '349.44'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.84'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.11%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.999'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.526'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.68%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '63.85'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.92%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.998'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.20'
$ws.Range('D27').Style = "Normal"
$ws.Range('E28').Value = '  +9.51%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.0₃0802'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.79'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +6.76%  '
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '163.28'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.65%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '19.88'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.60%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.62'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +7.57%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.05'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +3.18%  '
$ws.Range('E36').Value = '  +6.65%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '339.42'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +10.71%  '
$ws.Range('E38').Value = '  +2.23%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '4.07'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +4.76%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.903'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +6.55%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '38.42'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.36%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.19'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +2.64%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.622'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +2.25%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '20.29'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.53%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0249'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +2.73%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0563'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +2.21%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '132.78'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.63%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0994'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.83%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '20.54'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.24%  '
$ws.Range('E50').Value = '  +0.42%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.088.42'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +2.54%  '
